# "Generate Report for Handback" - update localization-status workbook
# with handback results for the two files (3f2c5e87..., c20530a5...) in
# both the zh-cn and de-de sheets, plus refresh the Overview status text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# 1. Overview sheet + per-language Status column: "Ready for handoff"
#    -> "Handed back: in sync with en-US" for both files/rows.
# ---------------------------------------------------------------
$ws1.Range("E2").Value = $handedBack
$ws1.Range("F2").Value = $handedBack
$ws1.Range("E3").Value = $handedBack
$ws1.Range("F3").Value = $handedBack

$ws2.Range("C2").Value = $handedBack
$ws2.Range("C3").Value = $handedBack

$ws3.Range("C2").Value = $handedBack
$ws3.Range("C3").Value = $handedBack

# ---------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I), Latest Handback
#    File (J) and Latest Handback DateTime (K) for both rows.
# ---------------------------------------------------------------
$ws2.Range("I2").Value = "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md"
$ws2.Range("I2").Font.Underline = 2
$ws2.Range("I2").Font.Color = 15570276
$ws2.Range("J2").Value = "3f2c5e87-344b-4b17-a884-2bd3de23af6e.10ff4ee0e067e1cb3467de8cdd4a0a01699cc5c2.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-10-17 16:46:53"

$ws2.Range("I3").Value = "c20530a5-4953-4e97-a05b-d115b33a7c19.md"
$ws2.Range("I3").Font.Underline = 2
$ws2.Range("I3").Font.Color = 15570276
$ws2.Range("J3").Value = "c20530a5-4953-4e97-a05b-d115b33a7c19.d384cac4cc4fdbc73c81afff80d1ed8ce6db92ed.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-10-17 16:46:53"

# ---------------------------------------------------------------
# 3. de-de sheet: same as above, different handback timestamp.
# ---------------------------------------------------------------
$ws3.Range("I2").Value = "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md"
$ws3.Range("I2").Font.Underline = 2
$ws3.Range("I2").Font.Color = 15570276
$ws3.Range("J2").Value = "3f2c5e87-344b-4b17-a884-2bd3de23af6e.10ff4ee0e067e1cb3467de8cdd4a0a01699cc5c2.de-de.xlf"
$ws3.Range("K2").Value = "2016-10-17 16:47:31"

$ws3.Range("I3").Value = "c20530a5-4953-4e97-a05b-d115b33a7c19.md"
$ws3.Range("I3").Font.Underline = 2
$ws3.Range("I3").Font.Color = 15570276
$ws3.Range("J3").Value = "c20530a5-4953-4e97-a05b-d115b33a7c19.d384cac4cc4fdbc73c81afff80d1ed8ce6db92ed.de-de.xlf"
$ws3.Range("K3").Value = "2016-10-17 16:47:31"

# ---------------------------------------------------------------
# 4. Rebuild the hyperlink collections so that the "Latest Target
#    File" cells (I2/I3) link to the same md files as the
#    "Source File Name" cells (A2/A3), in document order.
# ---------------------------------------------------------------
$url3f2c = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69bf272448e91a7d677ad7653765a29f0bbcb426/e2e/3f2c5e87-344b-4b17-a884-2bd3de23af6e.md"
$urlc205 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69bf272448e91a7d677ad7653765a29f0bbcb426/e2e/c20530a5-4953-4e97-a05b-d115b33a7c19.md"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $url3f2c, [Type]::Missing, [Type]::Missing, "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $url3f2c, [Type]::Missing, [Type]::Missing, "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlc205, [Type]::Missing, [Type]::Missing, "c20530a5-4953-4e97-a05b-d115b33a7c19.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlc205, [Type]::Missing, [Type]::Missing, "c20530a5-4953-4e97-a05b-d115b33a7c19.md")

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $url3f2c, [Type]::Missing, [Type]::Missing, "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $url3f2c, [Type]::Missing, [Type]::Missing, "3f2c5e87-344b-4b17-a884-2bd3de23af6e.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlc205, [Type]::Missing, [Type]::Missing, "c20530a5-4953-4e97-a05b-d115b33a7c19.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlc205, [Type]::Missing, [Type]::Missing, "c20530a5-4953-4e97-a05b-d115b33a7c19.md")

# ---------------------------------------------------------------
# 5. Column widths: the Status columns grew to fit the new, longer
#    status text, and the Target/Handback File columns grew to fit
#    the newly populated file names (capped at the report's usual
#    40-character-wide column cap).
# ---------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.17
$ws1.Columns.Item(6).ColumnWidth = 29.17

$ws2.Columns.Item(3).ColumnWidth = 29.17
$ws2.Columns.Item(9).ColumnWidth = 39.17
$ws2.Columns.Item(10).ColumnWidth = 39.17

$ws3.Columns.Item(3).ColumnWidth = 29.17
$ws3.Columns.Item(9).ColumnWidth = 39.17
$ws3.Columns.Item(10).ColumnWidth = 39.17

Write-Host "Handback report generated."
